# Updated cryptos list with GitHub Actions
#
# Applies refreshed price/volume figures (and two re-ranked coins whose
# name/link/price moved rows) to the crypto tracking sheet. Column D
# ("Price") holds price text that can look numeric (e.g. "247.71",
# "30.595.93"); it is written through a text-formatted round trip so
# Excel doesn't silently convert it to a number, then the cell style is
# restored to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D="30.595.93"; E="  +0.42%  " }
    @{ Row=3; D="1.923.41"; E="  -0.51%  " }
    @{ Row=4; E="  +0.07%  " }
    @{ Row=5; D="247.71"; E="  +2.89%  " }
    @{ Row=6; E="  +0.06%  " }
    @{ Row=7; D="0.4734"; E="  -0.36%  " }
    @{ Row=8; D="0.2921"; E="  +1.52%  " }
    @{ Row=9; D="0.06830"; E="  +2.37%  " }
    @{ Row=10; D="105.54"; E="  -1.88%  " }
    @{ Row=11; D="18.43"; E="  -3.67%  " }
    @{ Row=12; D="1.930.97"; E="  -0.10%  " }
    @{ Row=13; D="0.07726"; E="  +1.42%  " }
    @{ Row=14; D="5.323"; E="  +3.11%  " }
    @{ Row=15; D="0.6716"; E="  +1.23%  " }
    @{ Row=16; D="290.29"; E="  -5.51%  " }
    @{ Row=17; D="30.619.47"; E="  +0.44%  " }
    @{ Row=18; D="0.000007643"; E="  +0.47%  " }
    @{ Row=19; E="  +0.09%  " }
    @{ Row=20; D="12.97"; E="  -0.39%  " }
    @{ Row=21; D="5.562"; E="  +5.03%  " }
    @{ Row=22; D="2.176.19"; E="  -0.07%  " }
    @{ Row=23; D="1.000"; E="  -0.02%  " }
    @{ Row=24; D="6.486"; E="  +2.92%  " }
    @{ Row=25; D="9.538"; E="  +2.25%  " }
    @{ Row=26; D="167.18"; E="  -0.44%  " }
    @{ Row=27; D="20.79"; E="  +2.61%  " }
    @{ Row=28; D="2.140"; E="  +4.37%  " }
    @{ Row=29; D="0.1070"; E="  -3.30%  " }
    @{ Row=30; E="  +2.66%  " }
    @{ Row=31; D="4.206"; E="  +2.53%  " }
    @{ Row=32; D="4.063"; E="  +3.14%  " }
    @{ Row=33; D="0.05035"; E="  +0.21%  " }
    @{ Row=34; D="0.7347"; E="  -1.24%  " }
    @{ Row=35; D="1.146"; E="  -0.91%  " }
    @{ Row=36; D="0.02059"; E="  +4.77%  " }
    @{ Row=37; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.742"; E="  -0.38%  " }
    @{ Row=38; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="0.9997"; E="  +0.03%  " }
    @{ Row=39; E="  -0.31%  " }
    @{ Row=40; D="111.90"; E="  +3.73%  " }
    @{ Row=41; D="2.048"; E="  +0.13%  " }
    @{ Row=42; D="0.4455"; E="  +6.03%  " }
    @{ Row=43; D="0.8723"; E="  -1.08%  " }
    @{ Row=44; D="5.897"; E="  +1.39%  " }
    @{ Row=46; D="67.88"; E="  -3.47%  " }
    @{ Row=47; D="7.289"; E="  +0.08%  " }
    @{ Row=48; D="9.363"; E="  +1.39%  " }
    @{ Row=49; D="0.1255"; E="  +3.46%  " }
    @{ Row=50; D="47.76"; E="  +12.89%  " }
    @{ Row=51; B="Decentraland"; C="https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"; D="0.4125"; E="  +6.99%  " }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($u.ContainsKey("B")) {
        $ws.Range("B$row").Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Range("C$row").Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Preserve this cell's existing "text" storage for the price value:
        # round-trip through a text number format so a numeric-looking
        # string (e.g. "247.71") isn't auto-coerced into a number, then
        # restore the default "Normal" style so no formatting residue
        # remains on the cell.
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E$row").Value = $u.E
    }
}
